# Add a new "2022-Q3" sheet before "2022-Q2" (all later quarter sheets shift right
# automatically) and insert its fund-holding detail data, then add the matching
# summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right before the current "2022-Q2" tab.
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

# Template sheet to copy the header/column-A cell styling from (any quarter sheet works).
$template = $wb.Worksheets.Item(3)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @(0, "161810", "银华内需精选混合（LOF）", "23.47", "94.62", "6.61", "1.5514", 7),
    @(1, "009394", "银华同力精选混合", "17.87", "94.50", "5.38", "0.9614", 7),
    @(2, "002207", "前海开源金银珠宝主题精选混合C", "6.72", "90.85", "7.65", "0.5141", 9),
    @(3, "001302", "前海开源金银珠宝主题精选混合A", "3.99", "90.85", "7.65", "0.3052", 9),
    @(4, "003304", "前海开源沪港深核心资源灵活配置混合A", "3.45", "90.59", "7.91", "0.2729", 4),
    @(5, "003305", "前海开源沪港深核心资源灵活配置混合C", "1.89", "90.59", "7.91", "0.1495", 4)
)

# Make B:G text-formatted so the numeric-looking strings (fund codes, percentages,
# market values, ...) are stored as text (matches how every other quarter sheet
# stores these columns - this also preserves leading zeros in fund codes).
$q3.Range("B2:G7").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Copy over the bold/bordered header style and the bold column-A style from the template
# sheet (freshly-added sheets start out unstyled).
for ($col = 2; $col -le 8; $col++) {
    $template.Cells.Item(1, $col).Copy()
    $q3.Cells.Item(1, $col).PasteSpecial(-4122)
}
for ($row = 2; $row -le 7; $row++) {
    $template.Cells.Item($row, 1).Copy()
    $q3.Cells.Item($row, 1).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3 and push
#    everything else down (the per-row index in column A also shifts by one).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 3.75
# The inserted row inherited header-row formatting for B:D; reset those to the plain
# (unstyled) look the other data rows use.
$summary.Range("B2:D2").Style = "Normal"
# ...but column A keeps the bold/bordered style used by every other row-index cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Renumber the row-index column (A) for the rows that shifted down one place.
for ($row = 3; $row -le 9; $row++) {
    $summary.Cells.Item($row, 1).Value = $row - 2
}
